$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4825.0835
$ws.Range("J64").Value = 4817.522
$ws.Range("L64").Value = 4817.522
$ws.Range("N64").Value = -5313.522

$ws.Range("H67").Value = 4825.0835
$ws.Range("J67").Value = 4817.522
$ws.Range("L67").Value = 4817.522
$ws.Range("N67").Value = -6533.522

$ws.Range("H112").Value = 6682.923
$ws.Range("J112").Value = 7441.5654
$ws.Range("L112").Value = 22324.6962
$ws.Range("N112").Value = -24540.6962

$ws.Range("H113").Value = 2960
$ws.Range("I113").Value = 2900
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2900
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 354
$ws.Range("N113").Value = -9508

$ws.Range("H129").Value = 1280.8064
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1280.8064
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3842.4192
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13842.4192

$ws.Range("H137").Value = 1887.7693
$ws.Range("I137").Value = 1827.6923
$ws.Range("J137").Value = 2007.9231
$ws.Range("K137").Value = 5483.0769
$ws.Range("L137").Value = 6023.7693
$ws.Range("M137").Value = -2933.0769
$ws.Range("N137").Value = -11123.7693

$ws.Range("H138").Value = 3036840
$ws.Range("I138").Value = 9095257
$ws.Range("J138").Value = 7631.386
$ws.Range("K138").Value = 27285771
$ws.Range("L138").Value = 22894.158
$ws.Range("M138").Value = -27280631
$ws.Range("N138").Value = -33174.158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 51580.785
$ws.Range("I32").Value = 46934.156
$ws.Range("J32").Value = 61390.332
$ws.Range("K32").Value = 46934.156
$ws.Range("L32").Value = 61390.332
$ws.Range("M32").Value = -46647.156
$ws.Range("N32").Value = -61964.332

$ws.Range("H61").Value = 2395.8667
$ws.Range("I61").Value = 1640.6086
$ws.Range("J61").Value = 4877.4287
$ws.Range("K61").Value = 1640.6086
$ws.Range("L61").Value = 4877.4287
$ws.Range("M61").Value = -1428.6086
$ws.Range("N61").Value = -5301.4287

$ws.Range("H74").Value = 1514.5
$ws.Range("I74").Value = 1451.2174
$ws.Range("J74").Value = 1999.6666
$ws.Range("K74").Value = 1451.2174
$ws.Range("L74").Value = 1999.6666
$ws.Range("M74").Value = -577.2174
$ws.Range("N74").Value = -3747.6666

$ws.Range("H77").Value = 1514.5
$ws.Range("I77").Value = 1451.2174
$ws.Range("J77").Value = 1999.6666
$ws.Range("K77").Value = 7256.087
$ws.Range("L77").Value = 9998.333000000001
$ws.Range("M77").Value = -2888.087
$ws.Range("N77").Value = -18734.333

$ws.Range("H123").Value = 30428
$ws.Range("J123").Value = 30428
$ws.Range("L123").Value = 30428
$ws.Range("N123").Value = -40228

$ws.Range("H136").Value = 2395.8667
$ws.Range("I136").Value = 1640.6086
$ws.Range("J136").Value = 4877.4287
$ws.Range("K136").Value = 4921.825800000001
$ws.Range("L136").Value = 14632.2861
$ws.Range("M136").Value = -2371.825800000001
$ws.Range("N136").Value = -19732.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3180.1133
$ws.Range("I134").Value = 1977.5264
$ws.Range("K134").Value = 5932.5792
$ws.Range("M134").Value = -3397.5792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 98.2
$ws.Range("I7").Value = 97
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 97
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = 16
$ws.Range("N7").Value = -326

$ws.Range("H31").Value = 6270.4
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 6270.4
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H58").Value = 1109.3182
$ws.Range("I58").Value = 1185.8438
$ws.Range("J58").Value = 905.25
$ws.Range("K58").Value = 1185.8438
$ws.Range("L58").Value = 905.25
$ws.Range("M58").Value = -982.8438000000001
$ws.Range("N58").Value = -1311.25

$ws.Range("H86").Value = 8910.462
$ws.Range("I86").Value = 17627.334
$ws.Range("J86").Value = 1438.8572
$ws.Range("K86").Value = 17627.334
$ws.Range("L86").Value = 1438.8572
$ws.Range("M86").Value = -16504.334
$ws.Range("N86").Value = -3684.8572

$ws.Range("H89").Value = 8910.462
$ws.Range("I89").Value = 17627.334
$ws.Range("J89").Value = 1438.8572
$ws.Range("K89").Value = 88136.67
$ws.Range("L89").Value = 7194.286
$ws.Range("M89").Value = -82520.67
$ws.Range("N89").Value = -18426.286

$ws.Range("H136").Value = 1109.3182
$ws.Range("I136").Value = 1185.8438
$ws.Range("J136").Value = 905.25
$ws.Range("K136").Value = 3557.5314
$ws.Range("L136").Value = 2715.75
$ws.Range("M136").Value = -1007.5314
$ws.Range("N136").Value = -7815.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 39.11111
$ws.Range("I38").Value = 38.857143
$ws.Range("J38").Value = 40
$ws.Range("K38").Value = 116.571429
$ws.Range("L38").Value = 120
$ws.Range("M38").Value = 230.428571
$ws.Range("N38").Value = -814

$ws.Range("H75").Value = 4730.091
$ws.Range("J75").Value = 7804.8335
$ws.Range("L75").Value = 23414.5005
$ws.Range("N75").Value = -25410.5005

$ws.Range("H78").Value = 4730.091
$ws.Range("J78").Value = 7804.8335
$ws.Range("L78").Value = 70243.5015
$ws.Range("N78").Value = -80227.5015

$ws.Range("H107").Value = 437.16
$ws.Range("I107").Value = 489.6842
$ws.Range("J107").Value = 270.83334
$ws.Range("K107").Value = 1469.0526
$ws.Range("L107").Value = 812.5000200000001
$ws.Range("M107").Value = 450.9474
$ws.Range("N107").Value = -4652.50002

$ws.Range("H113").Value = 4407.2188
$ws.Range("J113").Value = 4637.7
$ws.Range("L113").Value = 13913.1
$ws.Range("N113").Value = -18253.1

$ws.Range("H120").Value = 11406
$ws.Range("I120").Value = 9676.666999999999
$ws.Range("J120").Value = 14000
$ws.Range("K120").Value = 29030.001
$ws.Range("L120").Value = 42000
$ws.Range("M120").Value = -24192.001
$ws.Range("N120").Value = -51676

$ws.Range("H122").Value = 1637.55
$ws.Range("I122").Value = 606.1111
$ws.Range("J122").Value = 2481.4546
$ws.Range("K122").Value = 5454.9999
$ws.Range("L122").Value = 22333.0914
$ws.Range("M122").Value = -3004.9999
$ws.Range("N122").Value = -27233.0914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 7584.75
$ws.Range("I36").Value = 3332.6667
$ws.Range("K36").Value = 3332.6667
$ws.Range("M36").Value = -2847.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15829.777
$ws.Range("I40").Value = 17693.6
$ws.Range("J40").Value = 13500
$ws.Range("K40").Value = 17693.6
$ws.Range("L40").Value = 13500
$ws.Range("M40").Value = -17557.6
$ws.Range("N40").Value = -13772

$ws.Range("H122").Value = 13976728
$ws.Range("J122").Value = 5476.25
$ws.Range("L122").Value = 16428.75
$ws.Range("N122").Value = -21328.75

$ws.Range("H132").Value = 8399.526
$ws.Range("I132").Value = 11018.417
$ws.Range("J132").Value = 3910
$ws.Range("K132").Value = 33055.251
$ws.Range("L132").Value = 11730
$ws.Range("M132").Value = -30525.251
$ws.Range("N132").Value = -16790

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2040.9375
$ws.Range("I132").Value = 2244.7778
$ws.Range("J132").Value = 1778.8572
$ws.Range("K132").Value = 6734.3334
$ws.Range("L132").Value = 5336.571599999999
$ws.Range("M132").Value = -4204.3334
$ws.Range("N132").Value = -10396.5716
